$d = $word.ActiveDocument

# Locate the two relevant list-item paragraphs by a distinctive,
# language-neutral (digit-based) fragment of their text instead of a
# hard-coded paragraph index, so the script keeps working even if
# paragraph numbering shifts.
$idxOld = -1
$idxNew = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($idxOld -eq -1 -and $t.Contains("105657")) {
        $idxOld = $i
    }
    if ($idxNew -eq -1 -and $t.Contains("32958")) {
        $idxNew = $i
    }
}

# --- Change 1 -------------------------------------------------------
# Paragraph referencing Φ.353.1/324/105657/Δ1/8-10-2002: only the
# leading "Την με " wording changes to "Τη με " (rest of the
# paragraph / run structure stays untouched).
$p1 = $d.Paragraphs.Item($idxOld).Range
$p1.Find.Execute("Την με ", $true, $false, $false, $false, $false, $true, 1, $false, "Τη με ", 2)

# --- Change 2 -------------------------------------------------------
# Paragraph referencing the 2018 ministerial decision (Φ.350.2/…) is
# replaced wholesale by a new, single-run paragraph about the 2019
# decision (Φ.351.1/…), collapsing the old multi-run text.
$p2 = $d.Paragraphs.Item($idxNew).Range
$r2 = $d.Range($p2.Start, $p2.End - 1)
$r2.Text = "Τη με αριθ. Φ.351.1/11/48020/Ε3/28-3-2019 (ΑΔΑ: ΩΩΤΗ4653ΠΣ-ΒΔ3) Υπουργική Απόφαση με θέμα: «Τοποθέτηση Περιφερειακών Διευθυντών Εκπαίδευσης»"
